# Clears out the previously zero-filled forecast/placeholder columns that
# were left over from concatenating multiple balance sheets into a single
# worksheet. The affected cells become empty (no value) instead of holding
# literal 0s.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AY57:BK57").ClearContents()
$ws.Range("AY58:BK58").ClearContents()
$ws.Range("V64:AY64").ClearContents()
$ws.Range("AY71:BK71").ClearContents()
$ws.Range("AY72:BK72").ClearContents()
$ws.Range("AY73:BK73").ClearContents()
$ws.Range("AY77:BK77").ClearContents()
$ws.Range("AY78:BK78").ClearContents()
$ws.Range("AJ79:BK79").ClearContents()
